$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.869.97"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.230.50"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "268.70"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +5.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.11"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +13.16%  "
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +2.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.30"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +9.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0928"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.58"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +8.23%  "
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.564.22"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.251.53"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.795"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.864.27"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +2.54%  "
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.10"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.41"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +5.49%  "
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.90"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -4.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.61"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +18.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.88"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.55"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +6.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.61"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -4.08%  "
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.54"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0907"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("E34").Value = "  +3.72%  "
$ws.Range("E35").Value = "  +2.24%  "
$ws.Range("E36").Value = "  +3.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0357"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.37"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.35"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +18.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.54"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -4.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.61"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +7.47%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.38"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0994"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.31"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.23"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +8.09%  "
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("E50").Value = "  -9.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.52"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +3.71%  "
